# Apply the RW01/RW02/RW03/RW04/RW05/RWxx probate forms test-case changes
# to the "ScenarioMapping" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 26-43: rewrite the FeatureFileName (B) / ScenarioName (C) columns
# (ID in column A and SmokeTest/RegressionTest in D/E keep their values,
# except E31/E32 which flip from "Yes" back to "No").

$rows = @(
    # row, B (feature file),              C (scenario)
    @(26, "probateFormsRW01.feature", "Open Estate"),
    @(27, "probateFormsRW01.feature", "Verify RW01 form"),
    @(28, "probateFormsRW01.feature", "Reset the RW01 form"),
    @(29, "probateFormsRW02.feature", "Open Estate"),
    @(30, "probateFormsRW02.feature", "Verify RW02 form"),
    @(31, "probateFormsRW02.feature", "Reset the RW02 form"),
    @(32, "probateFormsRW03.feature", "Open Estate"),
    @(33, "probateFormsRW03.feature", "Verify RW03 form"),
    @(34, "probateFormsRW03.feature", "Reset the RW03 form"),
    @(35, "probateFormsRW04.feature", "Open Estate"),
    @(36, "probateFormsRW04.feature", "Verify RW04 form"),
    @(37, "probateFormsRW04.feature", "Reset the RW04 form"),
    @(38, "probateFormsRW05.feature", "Open Estate"),
    @(39, "probateFormsRW05.feature", "Verify RW05 form"),
    @(40, "probateFormsRW05.feature", "Reset the RW05 form"),
    @(41, "probateFormsRWxx.feature", "Open Estate"),
    @(42, "probateFormsRWxx.feature", "Verify RWxx form"),
    @(43, "probateFormsRWxx.feature", "Reset the RWxx form")
)

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = "No"
    $ws.Cells.Item($r, 5).Value = "No"
}

# TC_026..TC_042 IDs for rows 27-43 (row 26 keeps TC_025)
$ids = @(
    @(27, "TC_026"), @(28, "TC_027"), @(29, "TC_028"), @(30, "TC_029"),
    @(31, "TC_030"), @(32, "TC_031"), @(33, "TC_032"), @(34, "TC_033"),
    @(35, "TC_034"), @(36, "TC_035"), @(37, "TC_036"), @(38, "TC_037"),
    @(39, "TC_038"), @(40, "TC_039"), @(41, "TC_040"), @(42, "TC_041"),
    @(43, "TC_042")
)
foreach ($id in $ids) {
    $ws.Cells.Item($id[0], 1).Value = $id[1]
}

# Rows 44-50: blank filler rows (column A only), matching rows 33-35 before
for ($r = 44; $r -le 50; $r++) {
    $ws.Cells.Item($r, 1).Value = ""
}

# Update the data validation list range to cover the extended rows
$ws.Range("D1:D43").Validation.Delete()
$ws.Range("D1:D43").Validation.Add(3, 1, 1, "Yes,No")
$ws.Range("E2:E43").Validation.Delete()
$ws.Range("E2:E43").Validation.Add(3, 1, 1, "Yes,No")

# Restore view state: top-left cell + active selection
$ws.Application.ActiveWindow.ScrollRow = 30
$ws.Range("E43").Select()
